# Add a "Save" column (H) to the s_vals sheet, matching the existing
# header/row styling conventions already used in columns B-G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1: new header cell "Save" - reuse the same cell formatting (bold,
# bordered, centered) already applied to the other header cells by
# copying G1's format onto H1 instead of constructing a brand-new style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# H2:H7: new data column, populated with 0 (plain numeric, unstyled -
# same as the other numeric data columns B-G).
$ws.Range("H2:H7").Value = 0
